$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new rows (3:5) for QSFP28 SR4 LP / Gen2 / Error Free variants
$ws.Rows("3:5").Insert()

# The old "best fit" helper column (B, full PN) is no longer the primary
# lookup column - keep its data but hide it, and widen/introduce the new
# human readable description column (C).
$ws.Columns.Item(2).Hidden = $true
$ws.Columns.Item(3).ColumnWidth = 100
$ws.Columns.Item(4).ColumnWidth = 11

# New rows column A values
$ws.Range("A3").Value = "QSFP28 SR4 LP"
$ws.Range("A4").Value = "QSFP28 SR4 Gen2"
$ws.Range("A5").Value = "QSFP28 SR4 Error Free"

# Column C (new "full part number" description column) for every row
$ws.Range("C2").Value = "FTLC9551"
$ws.Range("C3").Value = "FTLC9552"
$ws.Range("C4").Value = "FTLC9553"
$ws.Range("C5").Value = "FTLC9554"
$ws.Range("C6").Value = "FTLC9141"
$ws.Range("C7").Value = " FCBN425QE; FCBN425QB; FCBR425QB"
$ws.Range("C8").Value = " FCBN425QP"
$ws.Range("C9").Value = "FTLC8221"
$ws.Range("C10").Value = "FTLD10CX1"
$ws.Range("C11").Value = "FCBND10CD1;FCBRD10CD1"
$ws.Range("C12").Value = "FTLD10CE3;FTLD10CD3"
$ws.Range("C13").Value = "FCBND10CD3"
$ws.Range("C14").Value = "FCBND12CD1"
$ws.Range("C15").Value = "FCBND12CD3"
$ws.Range("C16").Value = "FTLD12CL3"
$ws.Range("C17").Value = "FCBN410QB1"
$ws.Range("C18").Value = "FCBN410QE2"
$ws.Range("C19").Value = "FCCG410QD3;FCBG410QD3;FCBN410QD3;FCCN410QD3;FCCG414QD3;FCBG414QD3;FCBN414QD3;FCCN414QD3"
$ws.Range("C20").Value = "FTL410QD1;FTL410QE1;FTL410QX1"
$ws.Range("C21").Value = "FTL410QD2;FTL410QE2;FTL410QX2"
$ws.Range("C22").Value = "FTL410QD3;FTL410QE3;FTL410QX3"
$ws.Range("C23").Value = "FTL410QD4;FTL410QE4;FTL410QX4"
$ws.Range("C24").Value = "FBOPD10SL1;FBOTD10FL1;FBOTD10SE1;FBOTD10SH1;FBOTD10SM1;FBOTD10SL1"
$ws.Range("C25").Value = "FBRTP08CL1C00-C2"
$ws.Range("C26").Value = "FBTTP08CL1C00-C3 "
$ws.Range("C27").Value = "FCBG110SD1"
$ws.Range("C28").Value = "FCBG110SD2;FCBN125SD1"
$ws.Range("C29").Value = "FCBG125SD1;FCCG125SD1;FCBN125SD1;FCCN125SD1"
$ws.Range("C30").Value = "FBOTD25SL"
$ws.Range("C31").Value = "FBOTD25FL2"
$ws.Range("C33").Value = "FBOTD25FL3"

# Match the author's final cursor position
$ws.Range("C11").Select() | Out-Null
